$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.163.31"
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").Value = "2.323.82"
$ws.Range("E3").Value = "  +1.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "303.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.53%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.72"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.11%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +2.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.17"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.99%  "
$ws.Range("E11").Value = "  -0.45%  "
$ws.Range("E12").Value = "  -0.73%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "17.77"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.50%  "
$ws.Range("E14").Value = "  +2.36%  "
$ws.Range("D15").Value = "2.684.34"
$ws.Range("E15").Value = "  +1.08%  "
$ws.Range("D16").Value = "2.316.38"
$ws.Range("E16").Value = "  +0.85%  "
$ws.Range("E17").Value = "  -1.06%  "
$ws.Range("D18").Value = "43.093.35"
$ws.Range("E18").Value = "  +0.41%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.97"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.95%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.23"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.27%  "
$ws.Range("D21").Value = "0.0₃0914"
$ws.Range("E21").Value = "  +0.84%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.19"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.73%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "240.59"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.81%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.17"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.68%  "
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.53"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "168.11"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.42%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "34.35"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.77%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.21"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.86%  "
$ws.Range("E31").Value = "  -6.09%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.95"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +8.95%  "
$ws.Range("E33").Value = "  +2.48%  "
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.73"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.04%  "
$ws.Range("E36").Value = "  -0.38%  "
$ws.Range("E37").Value = "  +1.21%  "
$ws.Range("E38").Value = "  +2.19%  "
$ws.Range("E39").Value = "  +0.15%  "
$ws.Range("E40").Value = "  -0.17%  "
$ws.Range("E41").Value = "  +0.36%  "
$ws.Range("D42").Value = "1.995.68"
$ws.Range("E42").Value = "  +0.29%  "
$ws.Range("E43").Value = "  +1.76%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.24"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.51%  "
$ws.Range("E45").Value = "  +1.32%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.61"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.35%  "
$ws.Range("E47").Value = "  +0.10%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "76.51"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +9.50%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "55.19"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.96%  "
$ws.Range("B50").Value = "HuobiToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.86"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +12.44%  "
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.549.77"
$ws.Range("E51").Value = "  +1.03%  "
